# Update the "想去人数" (want-to-go count) figures that were refreshed
# by the gh-pages data regeneration (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1994
$wsExpo.Range("F4").Value = 845
$wsExpo.Range("F5").Value = 1041

# Sheet "全部类型" (all types) - same events, shifted down one row
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1994
$wsAll.Range("F5").Value = 845
$wsAll.Range("F6").Value = 1041
